$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# reinterpreted as numbers (matches original inlineStr string cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.146.58'
$ws.Range("D3").Value = '3.402.09'
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '573.52'
$ws.Range("D6").Value = '142.55'
$ws.Range("D7").Value = '3.402.19'
$ws.Range("D8").Value = '1.00'
$ws.Range("D9").Value = '0.477'
$ws.Range("D10").Value = '7.60'
$ws.Range("D11").Value = '0.124'
$ws.Range("D12").Value = '0.397'
$ws.Range("D13").Value = '3.983.33'
$ws.Range("D14").Value = '0.125'
$ws.Range("D15").Value = '28.08'
$ws.Range("D16").Value = '0.0000172'
$ws.Range("D17").Value = '3.401.11'
$ws.Range("D18").Value = '61.170.23'
$ws.Range("D19").Value = '6.14'
$ws.Range("D20").Value = '13.89'
$ws.Range("D21").Value = '8.97'
$ws.Range("D22").Value = '383.86'
$ws.Range("D23").Value = '0.559'
$ws.Range("D24").Value = '74.54'
$ws.Range("D25").Value = '1.00'
$ws.Range("D26").Value = '0.0000118'
$ws.Range("D27").Value = '3.536.10'
$ws.Range("D28").Value = '0.180'
$ws.Range("D29").Value = '0.998'
$ws.Range("D30").Value = '7.39'
$ws.Range("D31").Value = '8.04'
$ws.Range("D32").Value = '2.17'
$ws.Range("D33").Value = '1.41'
$ws.Range("D34").Value = '0.999'
$ws.Range("D35").Value = '23.54'
$ws.Range("D36").Value = '7.03'
$ws.Range("D37").Value = '167.81'
$ws.Range("D38").Value = '3.433.34'
$ws.Range("D39").Value = '5.01'
$ws.Range("D40").Value = '1.49'
$ws.Range("D41").Value = '0.0775'
$ws.Range("D42").Value = '27.37'
$ws.Range("D43").Value = '0.783'
$ws.Range("D44").Value = '1.00'
$ws.Range("D45").Value = '4.45'
$ws.Range("D46").Value = '1.68'
$ws.Range("D47").Value = '1.14'
$ws.Range("D48").Value = '2.486.60'
$ws.Range("D49").Value = '6.84'
$ws.Range("D50").Value = '23.10'
$ws.Range("D51").Value = '0.0266'

# Reset style index back to default (no explicit style) now that the
# values are committed as text, to match the original unstyled cells.
$ws.Range("D2:D51").Style = "Normal"

# Column E (Volume) values are already text (percent strings with
# leading/trailing spaces), so a plain value assignment is sufficient.
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("E11").Value = '  -2.32%  '
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("E20").Value = '  -3.29%  '
$ws.Range("E21").Value = '  -5.00%  '
$ws.Range("E22").Value = '  -4.98%  '
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E26").Value = '  -5.07%  '
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").Value = '  -3.27%  '
$ws.Range("E31").Value = '  -2.93%  '
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  -3.00%  '
$ws.Range("E40").Value = '  -5.38%  '
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("E43").Value = '  -2.52%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("E46").Value = '  -3.84%  '
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("E48").Value = '  -5.25%  '
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("E51").Value = '  +1.14%  '
